$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-66) holds the "Förändrad" (last changed) date.
# Update its serial date value from 45188 (2023-09-19) to 45189 (2023-09-20).
$ws.Range("C2:C66").Value = 45189
